$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.412.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.308.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +5.60%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.307.24"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.885.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.437.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000164"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.310.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "424.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("E21").Value = "  -3.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.90%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.461.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.207"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.69%  "
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.57%  "
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("E39").Value = "  -3.55%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.855.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.98%  "
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.752"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0660"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "311.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.94%  "
$ws.Range("E51").Value = "  -0.23%  "